$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap C6 <-> C8 (Minjung <-> Minal)
$ws.Range("C6").Value = "Minal"
$ws.Range("C8").Value = "Minjung"

# Swap C17 <-> C18 (Sungwoo <-> Fionna)
$ws.Range("C17").Value = "Fionna"
$ws.Range("C18").Value = "Sungwoo"
